$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeDataBatch16")

# Update the username values in column E (rows 2-4) with new data
$ws.Range("E2").Value = "fire147"
$ws.Range("E3").Value = "water258"
$ws.Range("E4").Value = "soil369"
